# Auto-update draw results: append the 2025-09-30 Pick 3 draw as a new row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 14

# Force the new row's cells to text format first so Excel does not
# reinterpret date-like / numeric-like strings (e.g. "2025-09-30",
# "250930") as a real date serial or a number - the source data keeps
# everything as plain text, matching the rest of the sheet.
$ws.Range("A" + $newRow + ":E" + $newRow).NumberFormat = "@"

$ws.Cells.Item($newRow, 1).Value = "2025-09-30"
$ws.Cells.Item($newRow, 2).Value = "Pick 3"
$ws.Cells.Item($newRow, 3).Value = "250930"
$ws.Cells.Item($newRow, 4).Value = "0-1-3"
$ws.Cells.Item($newRow, 5).Value = "2025-09-30T21:37:59.014+04:00"
